# Earnings.xlsx update — refresh the quarterly/annual EPS layout and
# relabel the hidden helper column (B) to the new data-provider field
# names, matching the upstream "Add files via upload" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Row 2: the column offsets that the quote lookup uses (C2:F2
#    move from a 4..1 countdown to a flat -1; G2/"current" stays 0).
# ---------------------------------------------------------------
$ws.Range("C2").Value = -1
$ws.Range("D2").Value = -1
$ws.Range("E2").Value = -1
$ws.Range("F2").Value = -1
$ws.Range("G2").Value = 0

# ---------------------------------------------------------------
# 2. Column B (hidden) holds the field name passed into the XLL
#    lookup. Rename each to the new provider's field names.
# ---------------------------------------------------------------
$ws.Range("B3").Value = "Date"
$ws.Range("B4").Value = "Actual"
$ws.Range("B5").Value = "Estimate"
$ws.Range("B6").Value = "Surprise"
$ws.Range("B7").Value = "SurprisePercentage"
$ws.Range("B8").Value = "Actual"

# ---------------------------------------------------------------
# 3. Row 8 ("year ago" row) flips its quarter offset from +4 to -4
#    for every column.
# ---------------------------------------------------------------
$ws.Range("C8").Formula = '=_xll.xlquoteEarnings(Symbol,$B8,C$2-4)'
$ws.Range("D8").Formula = '=_xll.xlquoteEarnings(Symbol,$B8,D$2-4)'
$ws.Range("E8").Formula = '=_xll.xlquoteEarnings(Symbol,$B8,E$2-4)'

# ---------------------------------------------------------------
# 4. Columns F and G (rows 3-8) are re-entered as (legacy CSE /
#    dynamic) array formulas — same lookup, now array-entered.
#    Row 8 also carries the new -4 offset.
# ---------------------------------------------------------------
$ws.Range("F3").FormulaArray = '=_xll.xlquoteEarnings(Symbol,$B3,F$2)'
$ws.Range("G3").FormulaArray = '=_xll.xlquoteEarnings(Symbol,$B3,G$2)'

$ws.Range("F4").FormulaArray = '=_xll.xlquoteEarnings(Symbol,$B4,F$2)'
$ws.Range("G4").FormulaArray = '=_xll.xlquoteEarnings(Symbol,$B4,G$2)'

$ws.Range("F5").FormulaArray = '=_xll.xlquoteEarnings(Symbol,$B5,F$2)'
$ws.Range("G5").FormulaArray = '=_xll.xlquoteEarnings(Symbol,$B5,G$2)'

$ws.Range("F6").FormulaArray = '=_xll.xlquoteEarnings(Symbol,$B6,F$2)'
$ws.Range("G6").FormulaArray = '=_xll.xlquoteEarnings(Symbol,$B6,G$2)'

$ws.Range("F7").FormulaArray = '=_xll.xlquoteEarnings(Symbol,$B7,F$2)/100'
$ws.Range("G7").FormulaArray = '=_xll.xlquoteEarnings(Symbol,$B7,G$2)/100'

$ws.Range("F8").FormulaArray = '=_xll.xlquoteEarnings(Symbol,$B8,F$2-4)'
$ws.Range("G8").FormulaArray = '=_xll.xlquoteEarnings(Symbol,$B8,G$2-4)'
